$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# --- Experimental = true ---------------------------------------------------
# Typing the bare word true/false into a cell makes Excel store it as a
# Boolean. The source sheet stores these flags as plain text ("true"/
# "false"), matching the rest of the column, so we build the literal via a
# quoted formula and immediately convert it back to a static value with a
# Paste Special (values only) - this keeps the cell's genuine text type and
# its existing style/formatting untouched.
$wsMeta.Range("B7").Value = "=""true"""
$wsMeta.Range("B7").Copy() | Out-Null
$wsMeta.Range("B7").PasteSpecial(-4163) | Out-Null

# --- Date refreshed ----------------------------------------------------------
$wsMeta.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# --- Compositional = false --------------------------------------------------
$wsMeta.Range("B18").Value = "=""false"""
$wsMeta.Range("B18").Copy() | Out-Null
$wsMeta.Range("B18").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0
